# Update the "Forecast Comparison" sheet:
#  - Insert a new "Week_Start_Date" column after "Week" (new column B),
#    shifting ASIN / MyForecast / ... / is_holiday_week one column to the right.
#  - Normalize the Week labels from zero-padded (W01, W02, ...) to non
#    zero-padded (W1, W2, ...) for weeks 1-9.
#  - Populate the new Week_Start_Date column with the corresponding
#    calendar date (stored as text, matching the week labels) for every
#    data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before the current column B ("ASIN"), pushing the
# existing B:I columns to C:J.
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Map of row -> (week label, week start date text)
$weekData = [ordered]@{
    2  = @("W1",  "2025-01-05")
    3  = @("W2",  "2025-01-12")
    4  = @("W3",  "2025-01-19")
    5  = @("W4",  "2025-01-26")
    6  = @("W5",  "2025-02-02")
    7  = @("W6",  "2025-02-09")
    8  = @("W7",  "2025-02-16")
    9  = @("W8",  "2025-02-23")
    10 = @("W9",  "2025-03-02")
    11 = @("W10", "2025-03-09")
    12 = @("W11", "2025-03-16")
    13 = @("W12", "2025-03-23")
    14 = @("W13", "2025-03-30")
    15 = @("W14", "2025-04-06")
    16 = @("W15", "2025-04-13")
    17 = @("W16", "2025-04-20")
}

foreach ($row in $weekData.Keys) {
    $label = $weekData[$row][0]
    $startDate = $weekData[$row][1]

    $ws.Cells.Item($row, 1).Value = $label

    $cell = $ws.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $startDate
}
